# Scheduled-runner update: refresh market-price derived columns (H-N)
# on the Diabolos_Profits workbook. Values below mirror the latest
# Universalis pull; only numeric price/profit cells move, nothing else.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 399
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").Value = ""   # was -1762, now cleared
$ws.Range("H40").Value = 2444.2222
$ws.Range("J40").Value = 2512.25
$ws.Range("L40").Value = 2512.25
$ws.Range("N40").Value = -2862.25
$ws.Range("H94").Value = 7349.8
$ws.Range("I94").Value = 2250
$ws.Range("K94").Value = 2250
$ws.Range("M94").Value = -1799
$ws.Range("H100").Value = 1400
$ws.Range("I100").Value = 1400
$ws.Range("K100").Value = 1400
$ws.Range("M100").Value = -859
$ws.Range("H116").Value = 23383324
$ws.Range("I116").Value = 20920226
$ws.Range("K116").Value = 20920226
$ws.Range("M116").Value = -20916784
$ws.Range("H118").Value = 815.1429000000001
$ws.Range("I118").Value = 345.8
$ws.Range("K118").Value = 1037.4
$ws.Range("M118").Value = 619.5999999999999
$ws.Range("H129").Value = 718.6667
$ws.Range("I129").Value = 718.6667
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 2156.0001
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = 2843.9999
$ws.Range("N129").Value = ""   # was -14760.25, now cleared
$ws.Range("H137").Value = 2849.4
$ws.Range("I137").Value = 1993.5
$ws.Range("K137").Value = 5980.5
$ws.Range("M137").Value = -3430.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1325.4166
$ws.Range("I2").Value = 1355
$ws.Range("K2").Value = 1355
$ws.Range("M2").Value = -1242
$ws.Range("H23").Value = 13003
$ws.Range("J23").Value = 13003
$ws.Range("L23").Value = 13003
$ws.Range("N23").Value = -13521
$ws.Range("H32").Value = 4256.68
$ws.Range("I32").Value = 4305.619
$ws.Range("K32").Value = 4305.619
$ws.Range("M32").Value = -4018.619
$ws.Range("H45").Value = 2750
$ws.Range("H116").Value = 1325.4166
$ws.Range("I116").Value = 1355
$ws.Range("K116").Value = 1355
$ws.Range("M116").Value = 939
$ws.Range("H132").Value = 58826260
$ws.Range("I132").Value = 66669228
$ws.Range("K132").Value = 200007684
$ws.Range("M132").Value = -200005154

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1325.4166
$ws.Range("I3").Value = 1355
$ws.Range("K3").Value = 1355
$ws.Range("M3").Value = -1241
$ws.Range("H134").Value = 6252025.5
$ws.Range("J134").Value = 2711.1428
$ws.Range("L134").Value = 8133.428400000001
$ws.Range("N134").Value = -13203.4284

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = ""   # was -640, now cleared
$ws.Range("N14").Value = ""   # was -1340, now cleared
$ws.Range("H16").Value = 1292.8334
$ws.Range("I16").Value = 1251.5
$ws.Range("K16").Value = 1251.5
$ws.Range("M16").Value = -964.5
$ws.Range("H113").Value = 1292.8334
$ws.Range("I113").Value = 1251.5
$ws.Range("K113").Value = 1251.5
$ws.Range("M113").Value = 918.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 1235.6666
$ws.Range("I8").Value = 1235.6666
$ws.Range("K8").Value = 3706.9998
$ws.Range("M8").Value = -3567.9998
$ws.Range("H29").Value = 50
$ws.Range("I29").Value = 50
$ws.Range("K29").Value = 150
$ws.Range("M29").Value = 127
$ws.Range("H107").Value = 804.26666
$ws.Range("I107").Value = 969.2857
$ws.Range("J107").Value = 659.875
$ws.Range("K107").Value = 2907.8571
$ws.Range("L107").Value = 1979.625
$ws.Range("M107").Value = -987.8571000000002
$ws.Range("N107").Value = -5819.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 16430
$ws.Range("J32").Value = 16430
$ws.Range("L32").Value = 16430
$ws.Range("N32").Value = -17022
$ws.Range("H75").Value = 50001
$ws.Range("J75").Value = 50001
$ws.Range("L75").Value = 50001
$ws.Range("N75").Value = -51749
$ws.Range("H78").Value = 50001
$ws.Range("J78").Value = 50001
$ws.Range("L78").Value = 150003
$ws.Range("N78").Value = -158739
$ws.Range("H102").Value = 1593.0834
$ws.Range("I102").Value = 1374.125
$ws.Range("K102").Value = 1374.125
$ws.Range("M102").Value = 247.875
$ws.Range("H113").Value = 2330.25
$ws.Range("I113").Value = 2107
$ws.Range("K113").Value = 2107
$ws.Range("M113").Value = 63
$ws.Range("H122").Value = 2959.0715
$ws.Range("I122").Value = 2744.7
$ws.Range("J122").Value = 3495
$ws.Range("K122").Value = 8234.099999999999
$ws.Range("L122").Value = 10485
$ws.Range("M122").Value = -5784.099999999999
$ws.Range("N122").Value = -15385
$ws.Range("H132").Value = 3154.2068
$ws.Range("I132").Value = 2268.8
$ws.Range("J132").Value = 5121.778
$ws.Range("K132").Value = 6806.400000000001
$ws.Range("L132").Value = 15365.334
$ws.Range("M132").Value = -4276.400000000001
$ws.Range("N132").Value = -20425.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 24911.953
$ws.Range("I61").Value = 20157.867
$ws.Range("K61").Value = 20157.867
$ws.Range("M61").Value = -19955.867
$ws.Range("H113").Value = 24911.953
$ws.Range("I113").Value = 20157.867
$ws.Range("K113").Value = 20157.867
$ws.Range("M113").Value = -17987.867
$ws.Range("H132").Value = 3953.4614
$ws.Range("I132").Value = 3749.625
$ws.Range("K132").Value = 11248.875
$ws.Range("M132").Value = -8718.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 13583
$ws.Range("J15").Value = 13583
$ws.Range("L15").Value = 13583
$ws.Range("N15").Value = -14159
$ws.Range("H62").Value = 1705176.1
$ws.Range("I62").Value = 2385146.5
$ws.Range("K62").Value = 2385146.5
$ws.Range("M62").Value = -2384522.5
$ws.Range("H65").Value = 1705176.1
$ws.Range("I65").Value = 2385146.5
$ws.Range("K65").Value = 11925732.5
$ws.Range("M65").Value = -11922612.5
$ws.Range("H100").Value = 963.94446
$ws.Range("I100").Value = 842.38464
$ws.Range("K100").Value = 1684.76928
$ws.Range("M100").Value = -1143.76928
$ws.Range("H122").Value = 2809.0908
$ws.Range("I122").Value = 2809.0908
$ws.Range("K122").Value = 8427.2724
$ws.Range("M122").Value = -5977.2724
